# EMPANADA Data Files Key - "Created figures with day 5 tests"
#
# 1) The Day3 (rows 38-47) / Day4 (rows 48-57) "Original/New File Name" and
#    "Speed" cells were still carrying the old/inconsistent style (index 5 -
#    empty font, defaults to the workbook's base font). Re-format them to
#    the Arial style (index 2) used everywhere else on the sheet, cell by
#    cell, via PasteSpecial(formats-only) so values are untouched and no new
#    style/font entries get created.
# 2) Append two new Day5 rows (58 and 59) describing the new test videos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Cells (old style 5) that need to move to the normal Arial style (2).
$cellsToRestyle = @(
    "A38","B38","E38",
    "A39","B39","E39",
    "A40","E40",
    "B41","E41",
    "B42","E42",
    "B43","E43",
    "B44","E44",
    "B45","E45",
    "B46","E46",
    "B47","E47",
    "A48","B48",
    "A49","B49",
    "A50","B50",
    "A51","B51",
    "A52","B52",
    "A53","B53",
    "A54","B54",
    "A55","B55",
    "A56","B56",
    "A57","B57"
)

# Source cell already formatted with style index 2 (Arial / theme color).
$ws.Range("D38").Copy()
foreach ($addr in $cellsToRestyle) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

# New Day5 data rows.
$ws.Range("A58").Value = "N/A"
$ws.Range("B58").Value = "Day5-Earth-8mms.mov"
$ws.Range("C58").Value = 5
$ws.Range("D58").Value = "Earth"
$ws.Range("E58").Value = 8
$ws.Range("F58").Value = 2
$ws.Range("G58").Value = 13

$ws.Range("A59").Value = "N/A"
$ws.Range("B59").Value = "Day5-Earth-2mms.mov"
$ws.Range("C59").Value = 5
$ws.Range("D59").Value = "Earth"
$ws.Range("E59").Value = 2
$ws.Range("F59").Value = 3
$ws.Range("G59").Value = 50

# Row 58 and most of row 59 use the normal Arial style (2).
$ws.Range("D38").Copy()
$newRowCellsStyle2 = @("A58","B58","C58","D58","E58","F58","G58","A59","C59","D59","E59","F59","G59")
foreach ($addr in $newRowCellsStyle2) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

# B59 (the new file name) keeps the older default style (5), matching the
# author's edit.
$ws.Range("A41").Copy()
$ws.Range("B59").PasteSpecial($xlPasteFormats)
